# Updates the cryptocurrency prices/volumes table to reflect the
# latest scrape (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.492.43'
$ws.Range("E2").Value = '  -2.53%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.117.59'
$ws.Range("E3").Value = '  -2.90%  '

# Row 4
$ws.Range("E4").Value = '  +0.35%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.66'
$ws.Range("E5").Value = '  +0.43%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.50'
$ws.Range("E6").Value = '  -5.57%  '

# Row 7
$ws.Range("E7").Value = '  +0.08%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.113.42'
$ws.Range("E8").Value = '  -3.11%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.519'
$ws.Range("E9").Value = '  -2.12%  '

# Row 10
$ws.Range("E10").Value = '  -3.94%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.29'
$ws.Range("E11").Value = '  -4.82%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.469'
$ws.Range("E12").Value = '  -2.95%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("E13").Value = '  -3.42%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.24'
$ws.Range("E14").Value = '  -6.16%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.639.66'
$ws.Range("E15").Value = '  -2.20%  '

# Row 16
$ws.Range("E16").Value = '  +2.34%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.542.33'
$ws.Range("E17").Value = '  -2.47%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.121.79'
$ws.Range("E18").Value = '  -2.44%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.81'
$ws.Range("E19").Value = '  -4.33%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '473.16'
$ws.Range("E20").Value = '  -3.17%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.48'
$ws.Range("E21").Value = '  -3.47%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.709'
$ws.Range("E22").Value = '  -2.37%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.85'
$ws.Range("E23").Value = '  +0.08%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.51'
$ws.Range("E24").Value = '  -4.26%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.16'
$ws.Range("E25").Value = '  -2.44%  '

# Row 26
$ws.Range("E26").Value = '  +0.33%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.78'
$ws.Range("E27").Value = '  -6.51%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.43'
$ws.Range("E28").Value = '  -4.67%  '

# Row 29
$ws.Range("B29").Value = 'NEARProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.00'
$ws.Range("E29").Value = '  +1.69%  '

# Row 30
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.119'
$ws.Range("E30").Value = '  -9.56%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.06'
$ws.Range("E31").Value = '  -10.97%  '

# Row 32
$ws.Range("E32").Value = '  +0.33%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.66'
$ws.Range("E33").Value = '  -4.31%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.11'
$ws.Range("E34").Value = '  -4.22%  '

# Row 35
$ws.Range("E35").Value = '  -0.09%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0783'
$ws.Range("E36").Value = '  +4.72%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.92'
$ws.Range("E37").Value = '  -4.55%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.56'
$ws.Range("E38").Value = '  -4.41%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '454.88'
$ws.Range("E39").Value = '  -5.32%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.97'
$ws.Range("E40").Value = '  -9.91%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0391'
$ws.Range("E41").Value = '  -4.46%  '

# Row 42
$ws.Range("E42").Value = '  -7.00%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.28'
$ws.Range("E43").Value = '  -3.52%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.845.68'
$ws.Range("E44").Value = '  -3.58%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.28'
$ws.Range("E45").Value = '  -8.49%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.265'
$ws.Range("E46").Value = '  -6.62%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.41'
$ws.Range("E47").Value = '  +0.77%  '

# Row 48
$ws.Range("E48").Value = '  +0.00%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.11'
$ws.Range("E49").Value = '  -6.14%  '

# Row 50
$ws.Range("E50").Value = '  -3.41%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '118.81'
$ws.Range("E51").Value = '  -1.87%  '
